$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Nome Completo" column first (the order in which these cells
# were edited determines the order new strings are appended to the shared
# string table, so it is reproduced here to match the saved workbook).
$ws.Range("B4").Value = "Bulk User 31"
$ws.Range("B6").Value = "Bulk User 51"
$ws.Range("B10").Value = "Bulk User 91"
$ws.Range("B8").Value = "Bulk User 71"
$ws.Range("B3").Value = "Bulk User 22"
$ws.Range("B5").Value = "Bulk User 42"
$ws.Range("B7").Value = "Bulk User 62"
$ws.Range("B9").Value = "Bulk User 82"
$ws.Range("B11").Value = "Bulk User 102"
$ws.Range("B2").Value = "Bulk User 12"

# Update the "Nome Utilizador" (username) column, top to bottom.
$ws.Range("A2").Value = "testebulkuser1"
$ws.Range("A3").Value = "testebulkuser2"
$ws.Range("A4").Value = "testebulkuser3"
$ws.Range("A5").Value = "testebulkuser4"
$ws.Range("A6").Value = "testebulkuser5"
$ws.Range("A7").Value = "testebulkuser6"
$ws.Range("A8").Value = "testebulkuser7"
$ws.Range("A9").Value = "testebulkuser8"
$ws.Range("A10").Value = "testebulkuser9"
$ws.Range("A11").Value = "testebulkuser10"

# Update the "Email" column, top to bottom.
$ws.Range("C2").Value = "testebulkuser1@gmail.com"
$ws.Range("C3").Value = "testebulkuser2@gmail.com"
$ws.Range("C4").Value = "testebulkuser3@gmail.com"
$ws.Range("C5").Value = "testebulkuser4@gmail.com"
$ws.Range("C6").Value = "testebulkuser5@gmail.com"
$ws.Range("C7").Value = "testebulkuser6@gmail.com"
$ws.Range("C8").Value = "testebulkuser7@gmail.com"
$ws.Range("C9").Value = "testebulkuser8@gmail.com"
$ws.Range("C10").Value = "testebulkuser9@gmail.com"
$ws.Range("C11").Value = "testebulkuser10@gmail.com"

# Leave the saved selection on D12, as in the final workbook state.
$ws.Range("D12").Select()
